# Adds the third student's data row to the Bnumber list, matching the
# order the data was typed in (First Name, Last Name, Middle Name, B Number)
# and leaves the selection on the Middle Name cell of the new row, just
# like the author's session did.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Test_three"
$ws.Range("C4").Value = "Last Name Test 3"
$ws.Range("B4").Value = "Middle three"
$ws.Range("D4").Value = 3

$ws.Range("B4").Select() | Out-Null
